$wb = $excel.ActiveWorkbook

# The localization report is being refreshed for archive: the file
# "7af22dd3-acfc-4a21-ab3d-e7c6d1ac95c7.md" has moved from "Ready for
# handoff" to "In Translation" for both locales, and on the Overview
# sheet as well.

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B4").Value = "In Translation"
$overview.Range("C4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B4").Value = "In Translation"
